$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (J1, K1) using the same header style as the rest of row 1 ---
$ws.Range("J1").Value = "City"
$ws.Range("K1").Value = "PIN Code"
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1:K1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Add (empty) J/K cells for the existing data rows 2-7 ---
# A lone leading apostrophe forces a blank *text* cell (matching the
# empty inlineStr cells already used elsewhere in the sheet, e.g. I2:I6)
# rather than clearing the cell to a truly blank/no-value state.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 10).Value = "'"
    $ws.Cells.Item($r, 11).Value = "'"
}

# --- New data rows 8-10 ---

# Row 8
$ws.Cells.Item(8, 1).Value = "KULDEVI"
$ws.Cells.Item(8, 2).Value = "+918401018206, +919979073732"
$ws.Cells.Item(8, 3).Value = "Not Found"
$ws.Cells.Item(8, 4).Value = "Not Found"
$ws.Cells.Item(8, 5).Value = "Bus Stand, Yard Shop, Yard Shop No"
$ws.Cells.Item(8, 6).Value = "GOPALBHAI +91 84010 18206 NAKULBHAI +91 99790 73732 RFC KULDEVI FRUIT COMPAN @ld Fruit MMarkket Yard Shop No 57, TRADERS @Opps ! IMPOTER Bus Stand; Gondal 360371."
$ws.Cells.Item(8, 7).Value = "GOPALBHAI +91 84010 18206 NAKULBHAI +91 99790 73732 RFC KULDEVI FRUIT COMPAN @ld Fruit MMarkket Yard Shop No 57, TRADERS @Opps ! IMPOTER Bus Stand; Gondal 360371."
$ws.Cells.Item(8, 8).Value = "2025-05-09 15:48:49"
$ws.Cells.Item(8, 9).Value = "GOPALBHAI +91 84010 18206 (size: 1.24), NAKULBHAI +91 99790 73732 (size: 1.00), RFC (size: 3.61), KULDEVI (size: 10.00), FRUIT COMPAN (size: 3.61)"
$ws.Cells.Item(8, 10).Value = "Not Found"
$ws.Cells.Item(8, 11).Value = 360371

# Row 9
$ws.Cells.Item(9, 1).Value = "KULDEVI"
$ws.Cells.Item(9, 2).Value = "+918401018206, +919979073732"
$ws.Cells.Item(9, 3).Value = "Not Found"
$ws.Cells.Item(9, 4).Value = "Not Found"
$ws.Cells.Item(9, 5).Value = "Bus Stand, Yard Shop, Yard Shop No"
$ws.Cells.Item(9, 6).Value = "GOPALBHAI +91 84010 18206 NAKULBHAI +91 99790 73732 RFC KULDEVI FRUIT COMPAN @ld Fruit MMarkket Yard Shop No 57, TRADERS @Opps ! IMPOTER Bus Stand; Gondal 360371."
$ws.Cells.Item(9, 7).Value = "GOPALBHAI +91 84010 18206 NAKULBHAI +91 99790 73732 RFC KULDEVI FRUIT COMPAN @ld Fruit MMarkket Yard Shop No 57, TRADERS @Opps ! IMPOTER Bus Stand; Gondal 360371."
$ws.Cells.Item(9, 8).Value = "2025-05-09 15:48:56"
$ws.Cells.Item(9, 9).Value = "GOPALBHAI +91 84010 18206 (size: 1.24), NAKULBHAI +91 99790 73732 (size: 1.00), RFC (size: 3.61), KULDEVI (size: 10.00), FRUIT COMPAN (size: 3.61)"
$ws.Cells.Item(9, 10).Value = "Not Found"
$ws.Cells.Item(9, 11).Value = 360371

# Row 10
$ws.Cells.Item(10, 1).Value = "KULDEVI"
$ws.Cells.Item(10, 2).Value = "+918401018206, +919979073732"
$ws.Cells.Item(10, 3).Value = "Not Found"
$ws.Cells.Item(10, 4).Value = "Not Found"
$ws.Cells.Item(10, 5).Value = "Bus Stand, Yard Shop, Yard Shop No"
$ws.Cells.Item(10, 6).Value = "GOPALBHAI +91 84010 18206 NAKULBHAI +91 99790 73732 RFC KULDEVI FRUIT COMPAN @ld Fruit MMarkket Yard Shop No 57, TRADERS @Opps ! IMPOTER Bus Stand; Gondal 360371."
$ws.Cells.Item(10, 7).Value = "GOPALBHAI +91 84010 18206 NAKULBHAI +91 99790 73732 RFC KULDEVI FRUIT COMPAN @ld Fruit MMarkket Yard Shop No 57, TRADERS @Opps ! IMPOTER Bus Stand; Gondal 360371."
$ws.Cells.Item(10, 8).Value = "2025-05-13 22:38:09"
$ws.Cells.Item(10, 9).Value = "GOPALBHAI +91 84010 18206 (size: 1.24), NAKULBHAI +91 99790 73732 (size: 1.00), RFC (size: 3.61), KULDEVI (size: 10.00), FRUIT COMPAN (size: 3.61)"
$ws.Cells.Item(10, 10).Value = "Not Found"
# K10 holds the PIN code as TEXT (not a number) - force text with a leading apostrophe
$ws.Cells.Item(10, 11).Value = "'360371"

Write-Host "Business card rows 8-10 and City/PIN Code columns added"
